$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and some B/C (Coin/Link) values are plain text that can
# look numeric (e.g. "0.9990", "25.952.03"). Force text format on D-column
# cells before assigning so Excel does not reinterpret them as numbers and
# silently drop significant trailing zeros / renormalize the value.
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D15","D16","D17","D18","D19","D20","D21","D22","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D46","D47","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "25.952.03"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "1.743.96"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "247.83"
$ws.Range("E5").Value = "  +4.88%  "
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.5051"
$ws.Range("E7").Value = "  -4.40%  "
$ws.Range("D8").Value = "0.2742"
$ws.Range("E8").Value = "  -2.81%  "
$ws.Range("D9").Value = "0.06185"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("B10").Value = "TRON"
$ws.Range("C10").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D10").Value = "0.07257"
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.737.92"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "0.6538"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("E13").Value = "  -1.55%  "
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").Value = "77.68"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").Value = "0.9994"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "0.9989"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "25.968.73"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "11.83"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").Value = "0.000006834"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").Value = "1.971.95"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "4.370"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "5.399"
$ws.Range("E24").Value = "  +3.26%  "
$ws.Range("D25").Value = "136.67"
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("D26").Value = "1.494"
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").Value = "15.24"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").Value = "1.777"
$ws.Range("E28").Value = "  -1.61%  "
$ws.Range("D29").Value = "105.38"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("D30").Value = "3.900"
$ws.Range("E30").Value = "  +2.78%  "
$ws.Range("D31").Value = "0.08238"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").Value = "3.634"
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").Value = "0.04677"
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("D34").Value = "2.654"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("D35").Value = "0.9930"
$ws.Range("E35").Value = "  -1.42%  "
$ws.Range("D36").Value = "0.6186"
$ws.Range("E36").Value = "  -2.61%  "
$ws.Range("D37").Value = "2.755"
$ws.Range("E37").Value = "  +1.88%  "
$ws.Range("D38").Value = "0.01612"
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("D39").Value = "1.927"
$ws.Range("E39").Value = "  -2.16%  "
$ws.Range("D40").Value = "0.9990"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").Value = "100.04"
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("D42").Value = "0.3919"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "0.7588"
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("D44").Value = "5.007"
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("D46").Value = "6.294"
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("D47").Value = "55.54"
$ws.Range("E47").Value = "  +1.87%  "
$ws.Range("E48").Value = "  -1.60%  "
$ws.Range("D49").Value = "30.59"
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "7.575"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "0.3433"
$ws.Range("E51").Value = "  -1.24%  "
